$wb = $excel.ActiveWorkbook

# --- BATT_CONNECTOR sheet: mark the 3 parts as "NO" in the NEW column ---
$ws1 = $wb.Worksheets.Item("BATT_CONNECTOR")
$ws1.Range("B7").Value = "NO"
$ws1.Range("B8").Value = "NO"
$ws1.Range("B9").Value = "NO"

# --- _HISTORY sheet: log this production run as version 4 ---
$ws2 = $wb.Worksheets.Item("_HISTORY")
$ws2.Range("A7").Value = 4
$ws2.Range("B7").Value = 45195
$ws2.Range("B7").NumberFormat = "d-mmm"
$ws2.Range("C7").Value = "DGB"
$ws2.Range("D7").Value = "Se actualiza columna NEW para nueva produccion SETI"

# --- restore the view/selection state seen after these edits ---
$ws2.Range("B8").Select() | Out-Null
$ws1.Range("B10").Select() | Out-Null
